$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# A1 is a brand new column ("Applicant ID"); B1..N1 get renamed/relabeled text.
$ws.Range("A1").Value = "Applicant ID"
$ws.Range("B1").Value = "Service ID"
$ws.Range("C1").Value = "Service Type ID"
$ws.Range("D1").Value = " Application_Form_Payment_Status"
$ws.Range("E1").Value = "Date_Of_Inspection"
$ws.Range("F1").Value = "Current_Step"
$ws.Range("G1").Value = "MSE_Are_Documents_Verified"
$ws.Range("H1").Value = " Finance_Is_Application_Fee_Verified"
$ws.Range("I1").Value = "Finance_Is_Processing_Fee_Verified"
$ws.Range("J1").Value = "Finance_Is_Inspection_Fee_Verified"
$ws.Range("K1").Value = "Inspection Status"
$ws.Range("L1").Value = "Are Equipment And Monitoring Fees Verified"
$ws.Range("M1").Value = "Area_Officer_Approval"
$ws.Range("N1").Value = "Marine_Hod_Approval"

# --- Data row (row 2) ---
# A2 already holds the value that now represents Applicant ID (unchanged: 1)
$ws.Range("B2").Value = 1
# C2 used to hold the (incorrect) text date "29/29/08" with a date-style; it now
# holds a plain numeric value and must lose that number formatting.
$ws.Range("C2").ClearFormats()
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
# E2 becomes a genuine date (serial 44958 = 2023-02-01) formatted with the
# long-date custom format that used to live on C2.
$ws.Range("E2").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Range("E2").Value = 44958
$ws.Range("F2").Value = 2
# G2..N2 keep their existing values (0/1 flags) - nothing to change there.

# --- Column widths / bestFit layout ---
# (Values below are chosen so that, after this engine's internal
# char-width rounding, the resulting stored <col> width lands as close as
# possible to the authored bestFit widths.)
# A brand new column width for the inserted Applicant ID column.
$ws.Columns.Item(1).ColumnWidth = 9.666666666666666
# Former column D's width (service type id) now (slightly wider) lives under column C.
$ws.Columns.Item(3).ColumnWidth = 13.666666666666666
# Former column B's width (app payment status) now lives under column D.
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668
# Date column is now wider because of the long-date format / longer header text.
$ws.Columns.Item(5).ColumnWidth = 23.833333333333336
# Former column E's width (current step) now lives under column F.
$ws.Columns.Item(6).ColumnWidth = 10.5
# Column B ("Service ID") is short and no longer needs a custom/bestFit width;
# shrink it back down towards the sheet's standard column width.
$ws.Columns.Item(2).ColumnWidth = 8.833333333333332

# --- Selection ---
$ws.Range("E2").Select()
